# Row 28 and row 29 had their data swapped (the two sightings were
# re-ordered). Apply the per-cell changes directly (row-level .Value
# array read/write isn't reliable on multi-cell ranges in this host, so
# each changed cell is written individually with its literal target
# value, matching the OOXML diff cell-by-cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28 becomes what row 29 used to contain ---
$ws.Range("A28").Value = 112170166
$ws.Range("B28").Value = 56446
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 100049
$ws.Range("F28").Value = "Spillkråka"
$ws.Range("G28").Value = "Dryocopus martius"
$ws.Range("H28").Value = "(Linnaeus, 1758)"
$ws.Range("M28").ClearContents()
$ws.Range("Q28").Value = 693564
$ws.Range("R28").Value = 6551561
$ws.Range("Z28").Value = "14:08"
$ws.Range("AB28").Value = "14:08"
$ws.Range("AC28").Value = "Högstubbe från gran med stora hackade hål från födosök."
$ws.Range("AJ28").ClearContents()
$ws.Range("AK28").ClearContents()
$ws.Range("AL28").ClearContents()
$ws.Range("AO28").ClearContents()

# --- Row 29 becomes what row 28 used to contain ---
$ws.Range("A29").Value = 112170169
$ws.Range("B29").Value = 8367
$ws.Range("D29").Value = "LC"
$ws.Range("E29").Value = 106554
$ws.Range("F29").Value = "Björksplintborre"
$ws.Range("G29").Value = "Scolytus ratzeburgii"
$ws.Range("H29").Value = "Janson, 1856"
$ws.Range("M29").Value = "äldre gnagspår"
$ws.Range("Q29").Value = 693585
$ws.Range("R29").Value = 6551594
$ws.Range("Z29").Value = "14:00"
$ws.Range("AB29").Value = "14:00"
$ws.Range("AC29").ClearContents()
$ws.Range("AJ29").Value = "björkar"
$ws.Range("AK29").Value = "Betula"
$ws.Range("AL29").Value = "Björklåga med delar av barken kvar."
$ws.Range("AO29").Value = "Betula # Björklåga med delar av barken kvar."
